# Apply name updates: every name now carries the first letter of the
# person's last name appended (e.g. "Buse" -> "BuseE").
# This touches the shift-roster cells (B3:H5) and the under-40-hours
# summary table (B7:B17 + the C7:C17 hour totals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Weekly shift roster (rows 3-5, columns B..H = Mon..Sun) ---
$ws.Range("B3").Value = "BuseE/MartinM/UrmoO/RobertM/MadisK/HelenaM/"
$ws.Range("C3").Value = "DaniellS/UkuJ/UrmoO/MartinM/BuseE/Anne-MaiP/"
$ws.Range("D3").Value = "Anne-MaiP/TanelM/TaaviP/Sven-ErvinP/TriinuS/DaniellS/"
$ws.Range("E3").Value = "TanelM/KevinV/Anne-MaiP/DaniellS/JoelK/JoonasK/"
$ws.Range("F3").Value = "HelenaJ/HelenaM/UkuJ/JoelK/UrmoO/RasmusR/"
$ws.Range("G3").Value = "JoonasK/JoelK/FredK/Sven-ErvinP/RasmusR/HelenaM/"
$ws.Range("H3").Value = "RasmusR/TriinuS/Siim-KaarelK/Anne-MaiP/KevinV/JoonasK/"

$ws.Range("B4").Value = "FredK/Sven-ErvinP/TanelM/SiimL/JoelK/DenizG/"
$ws.Range("C4").Value = "RobertM/MadisK/DenizG/SiimL/RasmusR/TimoK/"
$ws.Range("D4").Value = "UkuJ/MadisK/UrmoO/FredK/SiimL/RasmusR/"
$ws.Range("E4").Value = "MadisK/TriinuS/SiimL/Siim-KaarelK/TaaviP/TimoK/"
$ws.Range("F4").Value = "DaniellS/Anne-MaiP/TaaviP/DenizG/BuseE/Siim-SanderS/"
$ws.Range("G4").Value = "Siim-SanderS/TanelM/DaniellS/RobertM/MadisK/UrmoO/"
$ws.Range("H4").Value = "Siim-SanderS/Sven-ErvinP/UkuJ/TanelM/FredK/RobertM/"

$ws.Range("B5").Value = "Siim-SanderS/TaaviP/Siim-KaarelK/"
$ws.Range("C5").Value = "JoonasK/HelenaM/KevinV/"
$ws.Range("D5").Value = "DenizG/Siim-SanderS/BuseE/"
$ws.Range("E5").Value = "RobertM/MartinM/FredK/"
$ws.Range("F5").Value = "KevinV/SiimL/TriinuS/"
$ws.Range("G5").Value = "MartinM/HelenaJ/TaaviP/"
$ws.Range("H5").Value = "BuseE/JoelK/"

# --- "People with under 40 hours" summary table (rows 7-17) ---
$ws.Range("B7").Value = "Sven-ErvinP"
$ws.Range("C7").Value = 32
$ws.Range("B8").Value = "JoonasK"
$ws.Range("C8").Value = 32
$ws.Range("B9").Value = "HelenaM"
$ws.Range("C9").Value = 32
$ws.Range("B10").Value = "TimoK"
$ws.Range("C10").Value = 16
$ws.Range("B11").Value = "DenizG"
$ws.Range("C11").Value = 32
$ws.Range("B12").Value = "MartinM"
$ws.Range("C12").Value = 32
$ws.Range("B13").Value = "TriinuS"
$ws.Range("C13").Value = 32
$ws.Range("B14").Value = "Siim-KaarelK"
$ws.Range("C14").Value = 24
$ws.Range("B15").Value = "UkuJ"
$ws.Range("C15").Value = 32
$ws.Range("B16").Value = "KevinV"
$ws.Range("C16").Value = 32
$ws.Range("B17").Value = "HelenaJ"
$ws.Range("C17").Value = 16

# Column contents grew longer (last-name initial appended), so re-fit the
# "bestFit" columns B..H to their new widest content, mirroring Excel's
# automatic column auto-fit behaviour after editing.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
